# Add the full (non-portable) file path to the data sheet.
#
# The original sheet used a shared string with a templated path
# ("$HOME/research_offsite/external/bauer-lab/guideseq/fastq_files")
# in every row of the path_to_files (N) column. Replace it with the
# concrete absolute path used on the author's machine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPath = "/home/tib163/research_offsite/external/bauer-lab/guideseq/fastq_files"

$ws.Range("N2").Value = $newPath
$ws.Range("N3").Value = $newPath
$ws.Range("N4").Value = $newPath
$ws.Range("N5").Value = $newPath

# Update the view: scroll so column D is the leftmost visible column,
# and move the active selection to N6.
$ws.Range("N6").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
